# agregue cliente el pela cta 50815, al archivo grupo_clientes.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PREMIUM row (B3): append new client account 50815 to the existing list
$ws.Range("B3").Value = "33.237.251.308.946.950.952.958.963.965.969.10165.20101.20110.20128.20148.20164.20236.20241.20246.20271.20284.20293.20309.20351.20380.20163.50815"

# Move the active selection to B4, matching the recorded cursor position after the edit
$ws.Range("B4").Select()
